$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by exactly one day,
# preserving the time-of-day fraction and the existing date/time cell style.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Match the saved view state: whole column A selected, scrolled back to the
# top of the sheet (clears the old topLeftCell/D66 active-cell scroll state).
$ws.Columns("A:A").Select() | Out-Null
